# Natmi following Dr Hou advice
# Update ligand/receptor-expressing cell counts (E, K: 1 -> 3) and the
# dependent total/specificity-weighted statistics in rows 2-7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = @{
        E = 3
        G = 98.912777
        H = 296.738331
        I = 0.8120825131376513
        J = 0.8120825131376513
        K = 3
        M = 4.237840333333334
        N = 12.713521
        O = 0.4475167411500002
        P = 0.4475167411500002
        Q = 419.1765558526057
        R = 3772.589002673451
        S = 0.3634205198242639
        T = 0.3634205198242639
    }
    3 = @{
        E = 3
        G = 98.912777
        H = 296.738331
        I = 0.8120825131376513
        J = 0.8120825131376513
        K = 3
        M = 5.231839666666667
        N = 15.695519
        O = 0.5524832588499998
        P = 0.5524832588499998
        Q = 517.4957902487544
        R = 4657.462112238789
        S = 0.4486619933133874
        T = 0.4486619933133874
    }
    4 = @{
        E = 3
        G = 17.04862266666667
        H = 51.14586800000001
        I = 0.1399706767982279
        J = 0.1399706767982279
        K = 3
        M = 4.237840333333334
        N = 12.713521
        O = 0.4475167411500002
        P = 0.4475167411500002
        Q = 72.24934076458091
        R = 650.2440668812282
        S = 0.06263922113730291
        T = 0.06263922113730291
    }
    5 = @{
        E = 3
        G = 17.04862266666667
        H = 51.14586800000001
        I = 0.1399706767982279
        J = 0.1399706767982279
        K = 3
        M = 5.231839666666667
        N = 15.695519
        O = 0.5524832588499998
        P = 0.5524832588499998
        Q = 89.19566032949913
        R = 802.7609429654921
        S = 0.07733145566092503
        T = 0.07733145566092503
    }
    6 = @{
        E = 3
        G = 5.839988000000001
        H = 17.519964
        I = 0.0479468100641207
        J = 0.04794681006412069
        K = 3
        M = 4.237840333333334
        N = 12.713521
        O = 0.4475167411500002
        P = 0.4475167411500002
        Q = 24.74893669258267
        R = 222.740430233244
        S = 0.02145700018843333
        T = 0.02145700018843332
    }
    7 = @{
        E = 3
        G = 5.839988000000001
        H = 17.519964
        I = 0.0479468100641207
        J = 0.04794681006412069
        K = 3
        M = 5.231839666666667
        N = 15.695519
        O = 0.5524832588499998
        P = 0.5524832588499998
        Q = 30.55388087125734
        R = 274.984927841316
        S = 0.02648980987568737
        T = 0.02648980987568737
    }
}

foreach ($row in $updates.Keys) {
    $rowValues = $updates[$row]
    foreach ($col in $rowValues.Keys) {
        $ws.Range("$col$row").Value = $rowValues[$col]
    }
}
